# Applies the "Change year using for index to 2013 so aligns with heating
# demand data" edit described by the commit/diff:
#   - rename sheets
#   - turn the JRC tool URL on sheet1 into a real hyperlink
#   - rework the "typical years" test-params block (sheet2): relabel, drop
#     the hard-coded capacity test value, add capacity/kWh columns, insert
#     a 2013 row (with its heating-demand-aligned numbers), and add an
#     explanatory (italic) footnote

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Rename the worksheets
# ---------------------------------------------------------------
$wsRegion = $wb.Worksheets.Item(1)
$wsRegion.Name = "variability by region"

$wsYears = $wb.Worksheets.Item(2)
$wsYears.Name = "typical years"

# ---------------------------------------------------------------
# 2. sheet1 "variability by region": make the JRC tool URL a real
#    hyperlink
# ---------------------------------------------------------------
$wsRegion.Hyperlinks.Add($wsRegion.Range("E4"), "https://re.jrc.ec.europa.eu/pvg_tools/en/tools.html")

# ---------------------------------------------------------------
# 3. sheet2 "typical years": rework the test-params block
# ---------------------------------------------------------------

# 3a. Heading text now mentions the API
$wsYears.Range("A1").Value = "Test params with API"

# 3b. The capacity row no longer carries a single hard-coded test
#     value (C5); "kWh" moves up to become the units label for the
#     yearly-production column header row (C6), taking on the bold
#     label style already used by B3/B4/B5/B7.
$wsYears.Range("C5").ClearContents()

$wsYears.Range("B3").Copy()
$wsYears.Range("C6").PasteSpecial(-4122)
$wsYears.Range("C6").Value = "kWh"

# 3c. Row 7 heading gains the two capacity-scenario columns
$wsYears.Range("C7").ClearFormats()
$wsYears.Range("C7").Value = "4kW of capacity"
$wsYears.Range("D7").Value = "3kW of capacity"

# 3d. New D-column (3kW scenario) production figures for the existing
#     2020 row
$wsYears.Range("D9").Value = 2271.1216800000002

# 3e. Insert the new 2013 row (heating demand profiles are for 2013)
#     ahead of the 2012 row, pushing 2012 and the Average row down
$wsYears.Rows.Item(11).Insert()

$wsYears.Range("B11").Value = 2013
$wsYears.Range("C15").Copy()
$wsYears.Range("C11").PasteSpecial(-4122)
$wsYears.Range("D11").Value = 2239.14921

# 3f. Explanatory footnote about why 2013 replaced 2020, in italics
$wsYears.Range("B17").Value = "Initially used 2020, but now using 2013 as heating demand profiles are for 2013"
$wsYears.Range("B17").Font.Italic = $true

# ---------------------------------------------------------------
# 4. Restore the on-screen selections
# ---------------------------------------------------------------
$wsRegion.Range("E4").Select()

$wsYears.Activate()
$wsYears.Range("D25").Select()
